# Add a new time-tracking entry ("Modellierung dokumentieren") as a new
# week-closing row in the Arbeitsmatrix sheet. This mirrors inserting a row
# at row 48 (pushing the old rows 48-52 down to 49-53), filling the new
# row 48 with the task data, and re-establishing the "week total" column
# (L/M) formatting + formulas that mark the end of a week block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# Insert a new row before row 48; this shifts rows 48:52 down to 49:53,
# and Excel auto-extends dependent ranges (data validation sqref, etc.)
$ws.Rows.Item(48).Insert()

# Seed the new row's cell formatting from the row above (the previous last
# entry of the "Systemarchitektur" week) so fonts/number-formats/borders
# for the data columns match the rest of the table.
$ws.Range("A47:M47").Copy()
$ws.Range("A48:M48").PasteSpecial(-4122)

# New task data
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Konzeptuelles Design"
$ws.Range("C48").Value = "[TASK]"
$ws.Range("D48").Value = "Systemarchitektur"
$ws.Range("E48").Value = "Modellierung dokumentieren"
$ws.Range("F48").Value = 44325
$ws.Range("G48").Value = 44338
$ws.Range("J48").Value = 0.5
$ws.Range("K48").Value = 0.75

# Hours-worked formula (same pattern used throughout column I)
$ws.Range("I48").Formula = "=ROUNDUP(((SUM(K48-J48)*24*60/60)/0.25),0)*0.25"

# Week-total formulas (column L / M), same pattern as the other week-ending
# rows (e.g. L39/M39, L30/M30, ...)
$ws.Range("L48").Formula = "=SUM(H41:I48)"
$ws.Range("M48").Formula = "=SUM(L48+16)"

# Re-apply the "week end" border/format to L48 (thin left + double bottom)
# matching the style used on the other week-total cells, then drop the
# number-format that got auto-inferred onto M48 from the formula so it goes
# back to the sheet's default (General) formatting like the other M cells.
$ws.Range("L39").Copy()
$ws.Range("L48").PasteSpecial(-4122)
$ws.Range("M48").ClearFormats()

# Restore the scroll position / selection similar to the saved workbook view
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("H56").Select()

$excel.CutCopyMode = 0
